$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" column (C) for all data rows (2-13): 46059 -> 46060
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# Rows 5,6,7,8,9,11,13 are reshuffled (values moved between rows).
# Apply the new A (Beteckning), B (Datum) and G (Area (ha)) values per row.

# Row 5
$ws.Cells.Item(5, 1).Value = "A 50934-2024"
$ws.Cells.Item(5, 2).Value = 45602
$ws.Cells.Item(5, 7).Value = 0.6

# Row 6
$ws.Cells.Item(6, 1).Value = "A 46779-2025"
$ws.Cells.Item(6, 2).Value = 45926
$ws.Cells.Item(6, 7).Value = 1.5

# Row 7
$ws.Cells.Item(7, 1).Value = "A 31120-2023"
$ws.Cells.Item(7, 2).Value = 45113
$ws.Cells.Item(7, 7).Value = 0.2

# Row 8
$ws.Cells.Item(8, 1).Value = "A 56948-2025"
$ws.Cells.Item(8, 2).Value = 45978.64356481482
$ws.Cells.Item(8, 7).Value = 4.7

# Row 9
$ws.Cells.Item(9, 1).Value = "A 64431-2023"
$ws.Cells.Item(9, 2).Value = 45280
$ws.Cells.Item(9, 7).Value = 0.5

# Row 11
$ws.Cells.Item(11, 1).Value = "A 27724-2022"
$ws.Cells.Item(11, 2).Value = 44743.48386574074
$ws.Cells.Item(11, 7).Value = 1.3

# Row 13
$ws.Cells.Item(13, 1).Value = "A 64445-2023"
$ws.Cells.Item(13, 2).Value = 45280
$ws.Cells.Item(13, 7).Value = 3.7
